# Apply "repull data, push all data, mean calculation" update
# to column F (dSF) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -3
    4  = 5
    7  = -4
    8  = -2
    9  = 4
    10 = -4
    11 = 1
    12 = 2
    13 = -4
    14 = 1
    15 = -3
    16 = 5
    17 = 9
    18 = -1
    19 = -1
    21 = -3
    22 = 3
    23 = -6
    24 = 2
    25 = -2
    26 = -3
    27 = 3
    31 = -2
    33 = -3
    34 = 1
    35 = 3
    36 = -3
    37 = 3
    38 = 2
    39 = 1
    40 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
